{"js": "// Office.js (Word JavaScript API) edit script.\n// Body of: async (context) => { ... }\n\n// Mapping of old multiplication-problem text to new text, as found in the\n// diff. Each of these strings is unique in the document, so a direct\n// search/replace per pair is safe and unambiguous.\nconst replacements = [\n  ['92\u00d777=7084', '89\u00d760=5340'],\n  ['49\u00d719=931', '40\u00d726=1040'],\n  ['69\u00d745=3105', '86\u00d796=8256'],\n  ['72\u00d728=2016', '54\u00d739=2106'],\n  ['18\u00d797=1746', '78\u00d764=4992'],\n  ['36\u00d772=2592', '34\u00d717=578'],\n  ['59\u00d757=3363', '50\u00d738=1900'],\n  ['35\u00d715=525', '69\u00d777=5313'],\n  ['28\u00d741=1148', '47\u00d719=893'],\n  ['86\u00d714=1204', '12\u00d749=588'],\n  ['56\u00d716=896', '74\u00d742=3108'],\n  ['82\u00d771=5822', '15\u00d769=1035'],\n  ['65\u00d779=5135', '93\u00d755=5115'],\n  ['97\u00d737=3589', '18\u00d728=504'],\n  ['39\u00d715=585', '50\u00d784=4200'],\n  ['33\u00d748=1584', '75\u00d778=5850'],\n  ['49\u00d730=1470', '36\u00d796=3456'],\n  ['48\u00d793=4464', '46\u00d721=966'],\n  ['64\u00d790=5760', '27\u00d755=1485'],\n  ['61\u00d784=5124', '15\u00d757=855'],\n  ['50\u00d714=700', '50\u00d786=4300'],\n  ['77\u00d772=5544', '13\u00d798=1274'],\n  ['63\u00d768=4284', '39\u00d765=2535'],\n  ['22\u00d755=1210', '16\u00d786=1376'],\n  ['74\u00d765=4810', '79\u00d737=2923'],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Word COM interop edit script.\n# $word is the Application, $d / $word.ActiveDocument is the Document.\n\n$d = $word.ActiveDocument\n\n# Mapping of old multiplication-problem text to new text, as found in the\n# diff. Each of these strings is unique within the document, so a simple\n# Find/Replace (wdReplaceAll, restricted to whole-document Content range)\n# per pair is safe and unambiguous.\n$replacements = @(\n  @('92\u00d777=7084', '89\u00d760=5340'),\n  @('49\u00d719=931', '40\u00d726=1040'),\n  @('69\u00d745=3105', '86\u00d796=8256'),\n  @('72\u00d728=2016', '54\u00d739=2106'),\n  @('18\u00d797=1746', '78\u00d764=4992'),\n  @('36\u00d772=2592', '34\u00d717=578'),\n  @('59\u00d757=3363', '50\u00d738=1900'),\n  @('35\u00d715=525', '69\u00d777=5313'),\n  @('28\u00d741=1148', '47\u00d719=893'),\n  @('86\u00d714=1204', '12\u00d749=588'),\n  @('56\u00d716=896', '74\u00d742=3108'),\n  @('82\u00d771=5822', '15\u00d769=1035'),\n  @('65\u00d779=5135', '93\u00d755=5115'),\n  @('97\u00d737=3589', '18\u00d728=504'),\n  @('39\u00d715=585', '50\u00d784=4200'),\n  @('33\u00d748=1584', '75\u00d778=5850'),\n  @('49\u00d730=1470', '36\u00d796=3456'),\n  @('48\u00d793=4464', '46\u00d721=966'),\n  @('64\u00d790=5760', '27\u00d755=1485'),\n  @('61\u00d784=5124', '15\u00d757=855'),\n  @('50\u00d714=700', '50\u00d786=4300'),\n  @('77\u00d772=5544', '13\u00d798=1274'),\n  @('63\u00d768=4284', '39\u00d765=2535'),\n  @('22\u00d755=1210', '16\u00d786=1376'),\n  @('74\u00d765=4810', '79\u00d737=2923')\n)\n\nforeach ($pair in $replacements) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n\n  $range = $d.Content\n  $find = $range.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $oldText\n  $find.Replacement.Text = $newText\n  $find.Forward = $true\n  $find.Wrap = 1  # wdFindContinue\n  $find.Format = $false\n  $find.MatchCase = $true\n  $find.MatchWholeWord = $false\n  $find.MatchWildcards = $false\n\n  $find.Execute($find.Text, $find.MatchCase, $find.MatchWholeWord, $find.MatchWildcards, $false, $false, $find.Forward, $find.Wrap, $false, $find.Replacement.Text, 2) | Out-Null\n}\n"}
